# edit.ps1
# Applies the "Adds pico container and DI" commit:
#   1) Nudges the "Objetivos" group box (slide 2) down slightly
#      (a:off y 1653233 -> 1662469 EMU).
#   2) Re-splits the run "Los miembros de la clase son los WebElements de
#      la pagina.Cada Pagina o seccion de la misma se implementa como una
#      clase" (slide 3, body placeholder, 2nd paragraph) into six runs
#      along the same word boundaries PowerPoint's own re-save produced,
#      without altering the rendered text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Group shape vertical offset (slide with "Objetivos" title).
# ---------------------------------------------------------------------
$targetLeftEmu = 1481816
$targetTopEmu  = 1653233
$newTopEmu     = 1662469

$groupShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.Type -eq 6) {
            $leftEmu = [Math]::Round($sh.Left * 12700)
            $topEmu  = [Math]::Round($sh.Top * 12700)
            if ($leftEmu -eq $targetLeftEmu -and $topEmu -eq $targetTopEmu) {
                $groupShape = $sh
            }
        }
    }
}

if ($groupShape -ne $null) {
    $groupShape.Top = $newTopEmu / 12700.0
}

# ---------------------------------------------------------------------
# 2) Re-split the "WebElements"/"Cada" paragraph on the Page Object
#    Model slide.
# ---------------------------------------------------------------------
$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTextFrame -ne 0) {
            if ($sh.TextFrame.HasText -ne 0) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t.Contains("WebElements")) {
                    $targetShape = $sh
                    $targetSlide = $sl
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $para = $tr.Paragraphs($pi)
        if ($para.Text.Contains("WebElements")) {
            # Re-assert the existing font size on each new sub-run to force
            # PowerPoint to split the run at that boundary while keeping
            # every inherited rPr attribute untouched.
            $para.Characters(1, 33).Font.Size  = 28   # "Los miembros de la clase son los "
            $para.Characters(34, 11).Font.Size = 28   # "WebElements"
            $para.Characters(45, 7).Font.Size  = 28   # " de la "
            $para.Characters(52, 7).Font.Size  = 28   # "pagina."
            $para.Characters(59, 4).Font.Size  = 28   # "Cada"
            $para.Characters(63, 58).Font.Size = 28   # " Pagina o seccion de la misma se implementa como una clase"
        }
    }
}
